$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (Förändrad) for rows 2 through 23 from 45175 to 45183
for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45175) {
        $cell.Value = 45183
    }
}
